$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (changed) date column for rows 2-7 from 45207 to 45208
foreach ($row in 2..7) {
    $ws.Cells.Item($row, 3).Value = 45208
}
